$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("output")

$ws.Range("A4").Value = "ddf--list--geo--global.csv"
$ws.Range("B4").Value = "latitude"
$ws.Range("C4").Value = "global"

$ws.Range("A5").Value = "ddf--list--geo--global.csv"
$ws.Range("B5").Value = "longitude"
$ws.Range("C5").Value = "global"

$lo = $ws.ListObjects.Item("list_index")
$lo.Resize($ws.Range("A1:C5"))
